$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Good Morning" (R10's greeting) is replaced with "GIT UPDATE".
$ws.Range("E8").Value = "GIT UPDATE"

# Leave the sheet with E8 selected, matching the saved view state.
$ws.Range("E8").Select()
